# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect the latest generated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    4  = 259
    6  = 205
    8  = 55
    14 = 2038
    16 = 17
    17 = 500
    18 = 473
    19 = 146
    20 = 72
    23 = 1533
    24 = 3734
    26 = 58
    28 = 1116
    29 = 102
    30 = 1901
    33 = 66
    35 = 406
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "全部类型" (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    4  = 259
    6  = 205
    8  = 55
    14 = 2038
    17 = 17
    18 = 500
    19 = 473
    20 = 146
    21 = 72
    24 = 1533
    25 = 3734
    27 = 58
    29 = 1116
    30 = 102
    31 = 1901
    34 = 66
    36 = 406
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}

$wb.Save()
